# Minor picture placement corrections (nudge a handful of picture/
# textbox offsets on slides 11-14). Shape positions in the PowerPoint
# object model are expressed in points, while the target offsets from
# the diff are given in EMU (1 pt = 12700 EMU), so each target EMU
# value below is converted to the closest points literal that survives
# the COM layer's internal float32 round-trip back to the exact EMU.

$p = $ppt.ActivePresentation

# Slide 11 - "Shape 358" (coPaperDBLPspeedup.png): x 228600 -> 237478
$s11 = $p.Slides.Item(11)
$sh358 = $s11.Shapes.Item(2)
$sh358.Left = 18.699055718110237

# Slide 12 - "Shape 366" (coPaperDBLPruntimes_matching.png): x 152400 -> 196790
$s12 = $p.Slides.Item(12)
$sh366 = $s12.Shapes.Item(2)
$sh366.Left = 15.49527599055118

# Slide 13 - "Shape 373" (coPaperDBLPspeedup_matching.png): x 200336 -> 235848, y 1702600 -> 1835770
$s13 = $p.Slides.Item(13)
$sh373 = $s13.Shapes.Item(1)
$sh373.Left = 18.570708661417324
$sh373.Top = 144.5488204976378

# Slide 13 - "Shape 376" (textbox): x 7014175 -> 7023053, y 640975 -> 649853
$sh376 = $s13.Shapes.Item(3)
$sh376.Left = 552.9963074125984
$sh376.Top = 51.16952755905512

# Slide 14 - "Shape 384" (textbox): x 7014175 -> 7023053, y 640975 -> 649853
$s14 = $p.Slides.Item(14)
$sh384 = $s14.Shapes.Item(2)
$sh384.Left = 552.9963074125984
$sh384.Top = 51.16952755905512

# Slide 14 - "Shape 385" (wikipediaspeedup_matching.png): x 228600 -> 246356, y 1690825 -> 1841751
$sh385 = $s14.Shapes.Item(3)
$sh385.Left = 19.39811043622047
$sh385.Top = 145.01976777952754
